$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation output values in row 2 (recomputed results)
$ws.Range("B2").Value = 19979382.17450935
$ws.Range("C2").Value = 617.23761494859673
$ws.Range("D2").Value = 0.96372302840079671
$ws.Range("E2").Value = 81493.378918601142
$ws.Range("F2").Value = 752.43916042029139
$ws.Range("G2").Value = 4.977148971513798
$ws.Range("H2").Value = 0.77410149382071514
$ws.Range("I2").Value = 0.79572239690916902
$ws.Range("J2").Value = 0.6693589797403684
$ws.Range("K2").Value = 886.68577143515631
$ws.Range("L2").Value = 1167.522837444762
$ws.Range("M2").Value = 1032.9487319997199
$ws.Range("N2").Value = 0.96372302840079671
$ws.Range("O2").Value = 0.96372302840027635
$ws.Range("P2").Value = 0.96372302839900503
$ws.Range("Q2").Value = 19982249.40916856
$ws.Range("R2").Value = 19981290.045560379
$ws.Range("S2").Value = 19980205.694626711
$ws.Range("T2").Value = 19979382.17450935
$ws.Range("U2").Value = 781999.63778310828
$ws.Range("V2").Value = 891474.28680362424
$ws.Range("W2").Value = 934505.23650035192
$ws.Range("X2").Value = 944800.79587955796
$ws.Range("Y2").Value = 4.977148971513798
$ws.Range("Z2").Value = 2.4124063996667511
$ws.Range("AA2").Value = 2.564742571856037
$ws.Range("AB2").Value = 110253.0074575888
$ws.Range("AC2").Value = 88288.016989655182
$ws.Range("AD2").Value = 81493.378918601142
$ws.Range("AE2").Value = 81493.378918601142
$ws.Range("AF2").Value = 896418.62872027443
$ws.Range("AG2").Value = 894425.10436861473
$ws.Range("AH2").Value = 877234.83445257659
$ws.Range("AI2").Value = 853289.10832661181

# Restore the last active selection saved with the workbook
$ws.Range("E13").Select()
